$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 225-226, pushing the existing data (old rows 225-319)
# down to new rows 227-321. Dimension grows from A1:T319 to A1:T321.
$ws.Rows("225:226").Insert()

# New row 225: Angeleno / Primera, 2023-03-24 (serial 45009)
$ws.Cells.Item(225, 1).Value = 4
$ws.Cells.Item(225, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(225, 3).Value = "Los Lagos"
$ws.Cells.Item(225, 4).Value = 45009
$ws.Cells.Item(225, 5).Value = 10
$ws.Cells.Item(225, 6).Value = "Fruta"
$ws.Cells.Item(225, 7).Value = 100103
$ws.Cells.Item(225, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(225, 9).Value = 100103002
$ws.Cells.Item(225, 10).Value = "Ciruela"
$ws.Cells.Item(225, 11).Value = "Angeleno"
$ws.Cells.Item(225, 12).Value = "Primera"
$ws.Cells.Item(225, 13).Value = 600
$ws.Cells.Item(225, 14).Value = 15000
$ws.Cells.Item(225, 15).Value = 16000
$ws.Cells.Item(225, 16).Value = 15500
$ws.Cells.Item(225, 17).Value = "`$/caja 14 kilos granel"
$ws.Cells.Item(225, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(225, 19).Value = 1107
$ws.Cells.Item(225, 20).Value = 14

# New row 226: Angeleno / Segunda, 2023-03-24 (serial 45009)
$ws.Cells.Item(226, 1).Value = 4
$ws.Cells.Item(226, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(226, 3).Value = "Los Lagos"
$ws.Cells.Item(226, 4).Value = 45009
$ws.Cells.Item(226, 5).Value = 10
$ws.Cells.Item(226, 6).Value = "Fruta"
$ws.Cells.Item(226, 7).Value = 100103
$ws.Cells.Item(226, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(226, 9).Value = 100103002
$ws.Cells.Item(226, 10).Value = "Ciruela"
$ws.Cells.Item(226, 11).Value = "Angeleno"
$ws.Cells.Item(226, 12).Value = "Segunda"
$ws.Cells.Item(226, 13).Value = 300
$ws.Cells.Item(226, 14).Value = 13000
$ws.Cells.Item(226, 15).Value = 13000
$ws.Cells.Item(226, 16).Value = 13000
$ws.Cells.Item(226, 17).Value = "`$/caja 14 kilos granel"
$ws.Cells.Item(226, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(226, 19).Value = 929
$ws.Cells.Item(226, 20).Value = 14
